# Auto-generated edit script.
# Source data in this sheet was re-exported/re-ordered upstream: a block of
# existing observation rows had their field-values permuted among each other
# (same rows, shuffled content) and three brand-new observation rows were
# appended at the bottom. This script reproduces the resulting cell values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 130746522
$ws.Cells.Item(4, 2).Value = 8451
$ws.Cells.Item(4, 4).Value = "LC"
$ws.Cells.Item(4, 5).Value = 106545
$ws.Cells.Item(4, 6).Value = "Mindre märgborre"
$ws.Cells.Item(4, 7).Value = "Tomicus minor"
$ws.Cells.Item(4, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(4, 10).Value = "'"
$ws.Cells.Item(4, 11).Value = "'"
$ws.Cells.Item(4, 12).Value = "'"
$ws.Cells.Item(4, 13).Value = "färska gnagspår"
$ws.Cells.Item(4, 14).Value = "'"
$ws.Cells.Item(4, 17).Value = 447866
$ws.Cells.Item(4, 18).Value = 6784597
$ws.Cells.Item(4, 32).Value = "'"
$ws.Cells.Item(5, 1).Value = 130746553
$ws.Cells.Item(5, 2).Value = 79243
$ws.Cells.Item(5, 4).Value = "NT"
$ws.Cells.Item(5, 5).Value = 6425
$ws.Cells.Item(5, 6).Value = "Garnlav"
$ws.Cells.Item(5, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(5, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 11).Value = ""
$ws.Cells.Item(5, 12).Value = ""
$ws.Cells.Item(5, 13).Value = ""
$ws.Cells.Item(5, 14).Value = ""
$ws.Cells.Item(5, 17).Value = 447903
$ws.Cells.Item(5, 18).Value = 6784473
$ws.Cells.Item(5, 32).Value = ""
$ws.Cells.Item(6, 1).Value = 130746558
$ws.Cells.Item(6, 2).Value = 79243
$ws.Cells.Item(6, 5).Value = 6425
$ws.Cells.Item(6, 6).Value = "Garnlav"
$ws.Cells.Item(6, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(6, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(6, 11).Value = ""
$ws.Cells.Item(6, 12).Value = ""
$ws.Cells.Item(6, 13).Value = ""
$ws.Cells.Item(6, 14).Value = ""
$ws.Cells.Item(6, 17).Value = 447718
$ws.Cells.Item(6, 18).Value = 6784468
$ws.Cells.Item(7, 1).Value = 130746493
$ws.Cells.Item(7, 2).Value = 57881
$ws.Cells.Item(7, 4).Value = "NT"
$ws.Cells.Item(7, 5).Value = 100049
$ws.Cells.Item(7, 6).Value = "Spillkråka"
$ws.Cells.Item(7, 7).Value = "Dryocopus martius"
$ws.Cells.Item(7, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(7, 10).Value = ""
$ws.Cells.Item(7, 13).Value = "äldre spår"
$ws.Cells.Item(7, 17).Value = 447766
$ws.Cells.Item(7, 18).Value = 6784433
$ws.Cells.Item(7, 32).Value = ""
$ws.Cells.Item(8, 1).Value = 130746520
$ws.Cells.Item(8, 2).Value = 8451
$ws.Cells.Item(8, 4).Value = "LC"
$ws.Cells.Item(8, 5).Value = 106545
$ws.Cells.Item(8, 6).Value = "Mindre märgborre"
$ws.Cells.Item(8, 7).Value = "Tomicus minor"
$ws.Cells.Item(8, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(8, 10).Value = "'"
$ws.Cells.Item(8, 11).Value = "'"
$ws.Cells.Item(8, 12).Value = "'"
$ws.Cells.Item(8, 13).Value = "äldre gnagspår"
$ws.Cells.Item(8, 14).Value = "'"
$ws.Cells.Item(8, 17).Value = 447846
$ws.Cells.Item(8, 18).Value = 6784643
$ws.Cells.Item(8, 32).Value = "'"
$ws.Cells.Item(10, 1).Value = 130746524
$ws.Cells.Item(10, 2).Value = 8451
$ws.Cells.Item(10, 4).Value = "LC"
$ws.Cells.Item(10, 5).Value = 106545
$ws.Cells.Item(10, 6).Value = "Mindre märgborre"
$ws.Cells.Item(10, 7).Value = "Tomicus minor"
$ws.Cells.Item(10, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(10, 10).Value = "'"
$ws.Cells.Item(10, 11).Value = "'"
$ws.Cells.Item(10, 12).Value = "'"
$ws.Cells.Item(10, 13).Value = "äldre gnagspår"
$ws.Cells.Item(10, 14).Value = "'"
$ws.Cells.Item(10, 17).Value = 447932
$ws.Cells.Item(10, 18).Value = 6784551
$ws.Cells.Item(10, 32).Value = "'"
$ws.Cells.Item(11, 1).Value = 130746530
$ws.Cells.Item(11, 17).Value = 447855
$ws.Cells.Item(11, 18).Value = 6784599
$ws.Cells.Item(12, 1).Value = 130746519
$ws.Cells.Item(12, 13).Value = "färska gnagspår"
$ws.Cells.Item(12, 17).Value = 447826
$ws.Cells.Item(12, 18).Value = 6784623
$ws.Cells.Item(13, 1).Value = 130746562
$ws.Cells.Item(13, 2).Value = 79243
$ws.Cells.Item(13, 4).Value = "NT"
$ws.Cells.Item(13, 5).Value = 6425
$ws.Cells.Item(13, 6).Value = "Garnlav"
$ws.Cells.Item(13, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(13, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(13, 10).Value = ""
$ws.Cells.Item(13, 11).Value = ""
$ws.Cells.Item(13, 12).Value = ""
$ws.Cells.Item(13, 13).Value = ""
$ws.Cells.Item(13, 14).Value = ""
$ws.Cells.Item(13, 17).Value = 447730
$ws.Cells.Item(13, 18).Value = 6784717
$ws.Cells.Item(13, 32).Value = ""
$ws.Cells.Item(15, 1).Value = 130746570
$ws.Cells.Item(15, 2).Value = 79243
$ws.Cells.Item(15, 4).Value = "NT"
$ws.Cells.Item(15, 5).Value = 6425
$ws.Cells.Item(15, 6).Value = "Garnlav"
$ws.Cells.Item(15, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(15, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(15, 10).Value = ""
$ws.Cells.Item(15, 11).Value = ""
$ws.Cells.Item(15, 12).Value = ""
$ws.Cells.Item(15, 13).Value = ""
$ws.Cells.Item(15, 14).Value = ""
$ws.Cells.Item(15, 17).Value = 447857
$ws.Cells.Item(15, 18).Value = 6784524
$ws.Cells.Item(15, 32).Value = ""
$ws.Cells.Item(16, 1).Value = 130746554
$ws.Cells.Item(16, 10).Value = "'"
$ws.Cells.Item(16, 11).Value = "'"
$ws.Cells.Item(16, 14).Value = "'"
$ws.Cells.Item(16, 17).Value = 447915
$ws.Cells.Item(16, 18).Value = 6784490
$ws.Cells.Item(16, 32).Value = "'"
$ws.Cells.Item(17, 1).Value = 130746500
$ws.Cells.Item(17, 2).Value = 57881
$ws.Cells.Item(17, 5).Value = 100049
$ws.Cells.Item(17, 6).Value = "Spillkråka"
$ws.Cells.Item(17, 7).Value = "Dryocopus martius"
$ws.Cells.Item(17, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(17, 10).Value = ""
$ws.Cells.Item(17, 12).Value = "'"
$ws.Cells.Item(17, 13).Value = "äldre spår"
$ws.Cells.Item(17, 17).Value = 447888
$ws.Cells.Item(17, 18).Value = 6784627
$ws.Cells.Item(17, 32).Value = ""
$ws.Cells.Item(19, 1).Value = 130746525
$ws.Cells.Item(19, 2).Value = 8451
$ws.Cells.Item(19, 4).Value = "LC"
$ws.Cells.Item(19, 5).Value = 106545
$ws.Cells.Item(19, 6).Value = "Mindre märgborre"
$ws.Cells.Item(19, 7).Value = "Tomicus minor"
$ws.Cells.Item(19, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(19, 10).Value = "'"
$ws.Cells.Item(19, 13).Value = "äldre gnagspår"
$ws.Cells.Item(19, 17).Value = 447933
$ws.Cells.Item(19, 18).Value = 6784552
$ws.Cells.Item(19, 32).Value = "'"
$ws.Cells.Item(23, 1).Value = 130746510
$ws.Cells.Item(23, 2).Value = 8451
$ws.Cells.Item(23, 4).Value = "LC"
$ws.Cells.Item(23, 5).Value = 106545
$ws.Cells.Item(23, 6).Value = "Mindre märgborre"
$ws.Cells.Item(23, 7).Value = "Tomicus minor"
$ws.Cells.Item(23, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(23, 10).Value = "'"
$ws.Cells.Item(23, 11).Value = "'"
$ws.Cells.Item(23, 12).Value = "'"
$ws.Cells.Item(23, 13).Value = "äldre gnagspår"
$ws.Cells.Item(23, 14).Value = "'"
$ws.Cells.Item(23, 17).Value = 447718
$ws.Cells.Item(23, 18).Value = 6784441
$ws.Cells.Item(23, 32).Value = "'"
$ws.Cells.Item(24, 1).Value = 130746556
$ws.Cells.Item(24, 2).Value = 79243
$ws.Cells.Item(24, 4).Value = "NT"
$ws.Cells.Item(24, 5).Value = 6425
$ws.Cells.Item(24, 6).Value = "Garnlav"
$ws.Cells.Item(24, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(24, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(24, 10).Value = ""
$ws.Cells.Item(24, 11).Value = ""
$ws.Cells.Item(24, 12).Value = ""
$ws.Cells.Item(24, 13).Value = ""
$ws.Cells.Item(24, 14).Value = ""
$ws.Cells.Item(24, 17).Value = 447876
$ws.Cells.Item(24, 18).Value = 6784529
$ws.Cells.Item(24, 32).Value = ""
$ws.Cells.Item(25, 1).Value = 130746564
$ws.Cells.Item(25, 2).Value = 79243
$ws.Cells.Item(25, 4).Value = "NT"
$ws.Cells.Item(25, 5).Value = 6425
$ws.Cells.Item(25, 6).Value = "Garnlav"
$ws.Cells.Item(25, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(25, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(25, 10).Value = ""
$ws.Cells.Item(25, 11).Value = ""
$ws.Cells.Item(25, 12).Value = ""
$ws.Cells.Item(25, 13).Value = ""
$ws.Cells.Item(25, 14).Value = ""
$ws.Cells.Item(25, 17).Value = 447866
$ws.Cells.Item(25, 18).Value = 6784648
$ws.Cells.Item(25, 32).Value = ""
$ws.Cells.Item(26, 1).Value = 130746569
$ws.Cells.Item(26, 2).Value = 79243
$ws.Cells.Item(26, 4).Value = "NT"
$ws.Cells.Item(26, 5).Value = 6425
$ws.Cells.Item(26, 6).Value = "Garnlav"
$ws.Cells.Item(26, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(26, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(26, 10).Value = ""
$ws.Cells.Item(26, 11).Value = ""
$ws.Cells.Item(26, 12).Value = ""
$ws.Cells.Item(26, 13).Value = ""
$ws.Cells.Item(26, 14).Value = ""
$ws.Cells.Item(26, 17).Value = 447856
$ws.Cells.Item(26, 18).Value = 6784518
$ws.Cells.Item(26, 32).Value = ""
$ws.Cells.Item(27, 1).Value = 130746515
$ws.Cells.Item(27, 2).Value = 8451
$ws.Cells.Item(27, 4).Value = "LC"
$ws.Cells.Item(27, 5).Value = 106545
$ws.Cells.Item(27, 6).Value = "Mindre märgborre"
$ws.Cells.Item(27, 7).Value = "Tomicus minor"
$ws.Cells.Item(27, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(27, 10).Value = "'"
$ws.Cells.Item(27, 11).Value = "'"
$ws.Cells.Item(27, 12).Value = "'"
$ws.Cells.Item(27, 13).Value = "färska gnagspår"
$ws.Cells.Item(27, 14).Value = "'"
$ws.Cells.Item(27, 17).Value = 447716
$ws.Cells.Item(27, 18).Value = 6784496
$ws.Cells.Item(27, 32).Value = "'"
$ws.Cells.Item(28, 1).Value = 130746518
$ws.Cells.Item(28, 2).Value = 8451
$ws.Cells.Item(28, 4).Value = "LC"
$ws.Cells.Item(28, 5).Value = 106545
$ws.Cells.Item(28, 6).Value = "Mindre märgborre"
$ws.Cells.Item(28, 7).Value = "Tomicus minor"
$ws.Cells.Item(28, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(28, 10).Value = "'"
$ws.Cells.Item(28, 11).Value = "'"
$ws.Cells.Item(28, 12).Value = "'"
$ws.Cells.Item(28, 13).Value = "äldre gnagspår"
$ws.Cells.Item(28, 14).Value = "'"
$ws.Cells.Item(28, 17).Value = 447815
$ws.Cells.Item(28, 18).Value = 6784612
$ws.Cells.Item(28, 32).Value = "'"
$ws.Cells.Item(36, 1).Value = 130746560
$ws.Cells.Item(36, 2).Value = 79243
$ws.Cells.Item(36, 4).Value = "NT"
$ws.Cells.Item(36, 5).Value = 6425
$ws.Cells.Item(36, 6).Value = "Garnlav"
$ws.Cells.Item(36, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(36, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(36, 10).Value = ""
$ws.Cells.Item(36, 11).Value = ""
$ws.Cells.Item(36, 12).Value = ""
$ws.Cells.Item(36, 13).Value = ""
$ws.Cells.Item(36, 14).Value = ""
$ws.Cells.Item(36, 17).Value = 447685
$ws.Cells.Item(36, 18).Value = 6784529
$ws.Cells.Item(36, 32).Value = ""
$ws.Cells.Item(37, 1).Value = 130746565
$ws.Cells.Item(37, 17).Value = 447912
$ws.Cells.Item(37, 18).Value = 6784599
$ws.Cells.Item(38, 1).Value = 130746561
$ws.Cells.Item(38, 17).Value = 447711
$ws.Cells.Item(38, 18).Value = 6784677
$ws.Cells.Item(39, 1).Value = 130746506
$ws.Cells.Item(39, 2).Value = 8451
$ws.Cells.Item(39, 4).Value = "LC"
$ws.Cells.Item(39, 5).Value = 106545
$ws.Cells.Item(39, 6).Value = "Mindre märgborre"
$ws.Cells.Item(39, 7).Value = "Tomicus minor"
$ws.Cells.Item(39, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(39, 10).Value = "'"
$ws.Cells.Item(39, 11).Value = "'"
$ws.Cells.Item(39, 12).Value = "'"
$ws.Cells.Item(39, 13).Value = "äldre gnagspår"
$ws.Cells.Item(39, 14).Value = "'"
$ws.Cells.Item(39, 17).Value = 447826
$ws.Cells.Item(39, 18).Value = 6784573
$ws.Cells.Item(39, 32).Value = "'"
$ws.Cells.Item(48, 1).Value = 131073028
$ws.Cells.Item(48, 2).Value = 8451
$ws.Cells.Item(48, 4).Value = "LC"
$ws.Cells.Item(48, 5).Value = 106545
$ws.Cells.Item(48, 6).Value = "Mindre märgborre"
$ws.Cells.Item(48, 7).Value = "Tomicus minor"
$ws.Cells.Item(48, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(48, 9).Value = "'"
$ws.Cells.Item(48, 13).Value = "färska gnagspår"
$ws.Cells.Item(48, 16).Value = "Evertsbergsvägen öst, Evertsbergsvägen öst, Dlr"
$ws.Cells.Item(48, 17).Value = 447573
$ws.Cells.Item(48, 18).Value = 6784308
$ws.Cells.Item(48, 19).Value = 10
$ws.Cells.Item(48, 20).Value = "Dalarna"
$ws.Cells.Item(48, 21).Value = "Älvdalen"
$ws.Cells.Item(48, 22).Value = "Dalarna"
$ws.Cells.Item(48, 23).Value = "Älvdalen"
$ws.Cells.Item(48, 25).Value = "'2026-02-08"
$ws.Cells.Item(48, 26).Value = "12:17"
$ws.Cells.Item(48, 27).Value = "'2026-02-08"
$ws.Cells.Item(48, 28).Value = "12:17"
$ws.Cells.Item(48, 30).Value = $false
$ws.Cells.Item(48, 31).Value = $false
$ws.Cells.Item(48, 33).Value = $false
$ws.Cells.Item(48, 46).Value = "'"
$ws.Cells.Item(48, 49).Value = "Eva Löfqvist"
$ws.Cells.Item(48, 50).Value = "Eva Löfqvist"
$ws.Cells.Item(48, 51).Value = "'"
$ws.Cells.Item(49, 1).Value = 131073156
$ws.Cells.Item(49, 2).Value = 79243
$ws.Cells.Item(49, 4).Value = "NT"
$ws.Cells.Item(49, 5).Value = 6425
$ws.Cells.Item(49, 6).Value = "Garnlav"
$ws.Cells.Item(49, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(49, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(49, 9).Value = "'"
$ws.Cells.Item(49, 16).Value = "Evertsbergsvägen öst, Evertsbergsvägen öst, Dlr"
$ws.Cells.Item(49, 17).Value = 447572
$ws.Cells.Item(49, 18).Value = 6784666
$ws.Cells.Item(49, 19).Value = 10
$ws.Cells.Item(49, 20).Value = "Dalarna"
$ws.Cells.Item(49, 21).Value = "Älvdalen"
$ws.Cells.Item(49, 22).Value = "Dalarna"
$ws.Cells.Item(49, 23).Value = "Älvdalen"
$ws.Cells.Item(49, 25).Value = "'2026-02-08"
$ws.Cells.Item(49, 26).Value = "12:24"
$ws.Cells.Item(49, 27).Value = "'2026-02-08"
$ws.Cells.Item(49, 28).Value = "12:24"
$ws.Cells.Item(49, 30).Value = $false
$ws.Cells.Item(49, 31).Value = $false
$ws.Cells.Item(49, 33).Value = $false
$ws.Cells.Item(49, 46).Value = "'"
$ws.Cells.Item(49, 49).Value = "Eva Löfqvist"
$ws.Cells.Item(49, 50).Value = "Eva Löfqvist"
$ws.Cells.Item(49, 51).Value = "'"
$ws.Cells.Item(50, 1).Value = 131073209
$ws.Cells.Item(50, 2).Value = 79243
$ws.Cells.Item(50, 4).Value = "NT"
$ws.Cells.Item(50, 5).Value = 6425
$ws.Cells.Item(50, 6).Value = "Garnlav"
$ws.Cells.Item(50, 7).Value = "Alectoria sarmentosa"
$ws.Cells.Item(50, 8).Value = "(Ach.) Ach."
$ws.Cells.Item(50, 9).Value = "'"
$ws.Cells.Item(50, 16).Value = "Evertsbergsvägen öst, Evertsbergsvägen öst, Dlr"
$ws.Cells.Item(50, 17).Value = 447552
$ws.Cells.Item(50, 18).Value = 6784648
$ws.Cells.Item(50, 19).Value = 10
$ws.Cells.Item(50, 20).Value = "Dalarna"
$ws.Cells.Item(50, 21).Value = "Älvdalen"
$ws.Cells.Item(50, 22).Value = "Dalarna"
$ws.Cells.Item(50, 23).Value = "Älvdalen"
$ws.Cells.Item(50, 25).Value = "'2026-02-08"
$ws.Cells.Item(50, 26).Value = "12:26"
$ws.Cells.Item(50, 27).Value = "'2026-02-08"
$ws.Cells.Item(50, 28).Value = "12:26"
$ws.Cells.Item(50, 30).Value = $false
$ws.Cells.Item(50, 31).Value = $false
$ws.Cells.Item(50, 33).Value = $false
$ws.Cells.Item(50, 46).Value = "'"
$ws.Cells.Item(50, 49).Value = "Eva Löfqvist"
$ws.Cells.Item(50, 50).Value = "Eva Löfqvist"
$ws.Cells.Item(50, 51).Value = "'"

